$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update shared string label: "Copper ores and concentrates" -> "Copper"
$ws1.Range("C4").Value = "Copper"
$ws2.Range("C4").Value = "Copper"

$arrS1R2 = New-Object 'object[,]' 1,90
$arrS1R2[0,0] = 1072.353295692879
$arrS1R2[0,1] = 1397.592550984388
$arrS1R2[0,2] = 1396.906498779885
$arrS1R2[0,3] = 1327.436354936901
$arrS1R2[0,4] = 1780.396702028428
$arrS1R2[0,5] = 1284.428452375686
$arrS1R2[0,6] = 2324.081358655799
$arrS1R2[0,7] = 2310.454967410925
$arrS1R2[0,8] = 2055.826773599356
$arrS1R2[0,9] = 1894.905341665539
$arrS1R2[0,10] = 4053.605641952215
$arrS1R2[0,11] = 2460.582577172025
$arrS1R2[0,12] = 4931.653758089969
$arrS1R2[0,13] = 5322.413908660931
$arrS1R2[0,14] = 5718.612607993894
$arrS1R2[0,15] = 12503.1699435628
$arrS1R2[0,16] = 13387.39326830922
$arrS1R2[0,17] = 14281.21778846244
$arrS1R2[0,18] = 15191.88837781273
$arrS1R2[0,19] = 16115.96600258726
$arrS1R2[0,20] = 16194.63490865154
$arrS1R2[0,21] = 17083.59080246898
$arrS1R2[0,22] = 17979.9011828671
$arrS1R2[0,23] = 18894.98023593718
$arrS1R2[0,24] = 19835.97117993856
$arrS1R2[0,25] = 9274.876520044618
$arrS1R2[0,26] = 9532.90480893649
$arrS1R2[0,27] = 9864.01774676778
$arrS1R2[0,28] = 10184.55480148741
$arrS1R2[0,29] = 10559.42529254623
$arrS1R2[0,30] = 6771.848656011638
$arrS1R2[0,31] = 6835.321037099697
$arrS1R2[0,32] = 7166.631508162252
$arrS1R2[0,33] = 7759.326183602462
$arrS1R2[0,34] = 8336.461818313506
$arrS1R2[0,35] = 11044.55847485559
$arrS1R2[0,36] = 12441.77413334231
$arrS1R2[0,37] = 14378.67869776967
$arrS1R2[0,38] = 16878.17027769432
$arrS1R2[0,39] = 18669.45809432919
$arrS1R2[0,40] = 22981.07668992534
$arrS1R2[0,41] = 23971.37023481244
$arrS1R2[0,42] = 24802.99732244884
$arrS1R2[0,43] = 25482.82038925307
$arrS1R2[0,44] = 26172.83913938083
$arrS1R2[0,45] = 27414.32540882629
$arrS1R2[0,46] = 27138.93639940287
$arrS1R2[0,47] = 25377.86381619268
$arrS1R2[0,48] = 22517.72673325493
$arrS1R2[0,49] = 21098.63861664585
$arrS1R2[0,50] = 20216.67603421636
$arrS1R2[0,51] = 20081.00806100281
$arrS1R2[0,52] = 19440.46143714973
$arrS1R2[0,53] = 18460.96432026254
$arrS1R2[0,54] = 18140.21304260692
$arrS1R2[0,55] = 17963.83647585458
$arrS1R2[0,56] = 18815.63982574294
$arrS1R2[0,57] = 20123.83979435458
$arrS1R2[0,58] = 21919.73256853574
$arrS1R2[0,59] = 23883.02251145838
$arrS1R2[0,60] = 25623.94081547318
$arrS1R2[0,61] = 27954.56778735401
$arrS1R2[0,62] = 30355.40261975654
$arrS1R2[0,63] = 32611.58445037794
$arrS1R2[0,64] = 34146.88689519095
$arrS1R2[0,65] = 34891.5727677519
$arrS1R2[0,66] = 35675.16333911721
$arrS1R2[0,67] = 36312.21076063113
$arrS1R2[0,68] = 36540.37801449133
$arrS1R2[0,69] = 35911.57011788043
$arrS1R2[0,70] = 34210.25787355839
$arrS1R2[0,71] = 32528.8442354809
$arrS1R2[0,72] = 31343.77141958112
$arrS1R2[0,73] = 30716.99926326832
$arrS1R2[0,74] = 30313.08139966867
$arrS1R2[0,75] = 29661.87144891542
$arrS1R2[0,76] = 29246.96719183742
$arrS1R2[0,77] = 29176.23747776269
$arrS1R2[0,78] = 29606.84253575559
$arrS1R2[0,79] = 30601.52982148424
$arrS1R2[0,80] = 31815.66021476475
$arrS1R2[0,81] = 33589.9244358727
$arrS1R2[0,82] = 35525.20814646035
$arrS1R2[0,83] = 37560.34104772157
$arrS1R2[0,84] = 39720.78912721969
$arrS1R2[0,85] = 41446.8227353573
$arrS1R2[0,86] = 43256.16286796908
$arrS1R2[0,87] = 44574.91062180322
$arrS1R2[0,88] = 45448.04477245767
$arrS1R2[0,89] = 46047.79029759407
$ws1.Range("D2:CO2").Value = $arrS1R2

$arrS1R3 = New-Object 'object[,]' 1,90
$arrS1R3[0,0] = 70.9098667458114
$arrS1R3[0,1] = 92.69577989651185
$arrS1R3[0,2] = 92.53143555583677
$arrS1R3[0,3] = 87.86447394665797
$arrS1R3[0,4] = 118.3849281223296
$arrS1R3[0,5] = 85.06207674038733
$arrS1R3[0,6] = 154.7765219514088
$arrS1R3[0,7] = 153.7091102740597
$arrS1R3[0,8] = 136.652725353146
$arrS1R3[0,9] = 125.7674442217332
$arrS1R3[0,10] = 270.6654012423808
$arrS1R3[0,11] = 163.8103713536561
$arrS1R3[0,12] = 329.6558393883638
$arrS1R3[0,13] = 355.8734260884138
$arrS1R3[0,14] = 382.4594554224848
$arrS1R3[0,15] = 837.956856854081
$arrS1R3[0,16] = 897.2973525566686
$arrS1R3[0,17] = 957.2876855605583
$arrS1R3[0,18] = 1018.418529479577
$arrS1R3[0,19] = 1080.456966172892
$arrS1R3[0,20] = 1085.726356490448
$arrS1R3[0,21] = 1145.433847300258
$arrS1R3[0,22] = 1205.638897725874
$arrS1R3[0,23] = 1267.114020856015
$arrS1R3[0,24] = 1330.342667918523
$arrS1R3[0,25] = 621.1844737261076
$arrS1R3[0,26] = 638.5344824115559
$arrS1R3[0,27] = 660.8308924368475
$arrS1R3[0,28] = 682.4076794276223
$arrS1R3[0,29] = 707.6589615660836
$arrS1R3[0,30] = 453.3399311568587
$arrS1R3[0,31] = 457.6076189312213
$arrS1R3[0,32] = 480.0066767806351
$arrS1R3[0,33] = 520.0947753887361
$arrS1R3[0,34] = 559.1188142578612
$arrS1R3[0,35] = 741.4269490232593
$arrS1R3[0,36] = 835.9285561235457
$arrS1R3[0,37] = 966.973190731153
$arrS1R3[0,38] = 1136.1147334345
$arrS1R3[0,39] = 1257.301627118642
$arrS1R3[0,40] = 1547.302395027694
$arrS1R3[0,41] = 1614.190103975942
$arrS1R3[0,42] = 1670.34750678821
$arrS1R3[0,43] = 1716.242226224635
$arrS1R3[0,44] = 1762.846447965382
$arrS1R3[0,45] = 1846.389023978494
$arrS1R3[0,46] = 1827.632793242692
$arrS1R3[0,47] = 1708.27971193741
$arrS1R3[0,48] = 1514.507753754608
$arrS1R3[0,49] = 1418.331047619856
$arrS1R3[0,50] = 1359.03755870426
$arrS1R3[0,51] = 1349.784623740363
$arrS1R3[0,52] = 1306.33215130483
$arrS1R3[0,53] = 1239.912604317023
$arrS1R3[0,54] = 1218.09091171905
$arrS1R3[0,55] = 1206.396780977021
$arrS1R3[0,56] = 1263.957050052245
$arrS1R3[0,57] = 1352.407483074634
$arrS1R3[0,58] = 1473.869612379449
$arrS1R3[0,59] = 1606.655716870921
$arrS1R3[0,60] = 1724.575517032822
$arrS1R3[0,61] = 1882.227145680058
$arrS1R3[0,62] = 2044.633212441249
$arrS1R3[0,63] = 2197.246131956082
$arrS1R3[0,64] = 2301.044457912811
$arrS1R3[0,65] = 2351.407454287733
$arrS1R3[0,66] = 2404.31135418801
$arrS1R3[0,67] = 2447.301378689442
$arrS1R3[0,68] = 2462.612576760546
$arrS1R3[0,69] = 2419.899330777072
$arrS1R3[0,70] = 2304.688258056013
$arrS1R3[0,71] = 2190.708606091831
$arrS1R3[0,72] = 2110.34784961533
$arrS1R3[0,73] = 2067.799982631145
$arrS1R3[0,74] = 2040.347487654176
$arrS1R3[0,75] = 1996.260681505848
$arrS1R3[0,76] = 1968.060547873602
$arrS1R3[0,77] = 1963.165155080213
$arrS1R3[0,78] = 1992.216015975052
$arrS1R3[0,79] = 2059.46107353669
$arrS1R3[0,80] = 2141.699809439252
$arrS1R3[0,81] = 2261.717592178501
$arrS1R3[0,82] = 2392.631977995993
$arrS1R3[0,83] = 2530.301941049948
$arrS1R3[0,84] = 2676.454125149572
$arrS1R3[0,85] = 2793.407254801385
$arrS1R3[0,86] = 2915.776968966183
$arrS1R3[0,87] = 3004.922752690742
$arrS1R3[0,88] = 3063.891805857525
$arrS1R3[0,89] = 3104.348952963342
$ws1.Range("D3:CO3").Value = $arrS1R3

$arrS1R4 = New-Object 'object[,]' 1,90
$arrS1R4[0,0] = 1624090.537329962
$arrS1R4[0,1] = 1458823.298962824
$arrS1R4[0,2] = 1273536.948864026
$arrS1R4[0,3] = 1477258.372532633
$arrS1R4[0,4] = 1623358.757813651
$arrS1R4[0,5] = 1542872.897228797
$arrS1R4[0,6] = 1677419.056944033
$arrS1R4[0,7] = 1538620.78706379
$arrS1R4[0,8] = 1528306.434827912
$arrS1R4[0,9] = 2038421.069231238
$arrS1R4[0,10] = 1934408.236466115
$arrS1R4[0,11] = 1916153.126861282
$arrS1R4[0,12] = 9931577.840286806
$arrS1R4[0,13] = 9982387.278428297
$arrS1R4[0,14] = 10022431.05349517
$arrS1R4[0,15] = 31279781.12871488
$arrS1R4[0,16] = 31344377.31664766
$arrS1R4[0,17] = 31419986.93323457
$arrS1R4[0,18] = 31505528.09695473
$arrS1R4[0,19] = 31599563.28161091
$arrS1R4[0,20] = 32549961.40126257
$arrS1R4[0,21] = 32660911.6093827
$arrS1R4[0,22] = 32769311.34836029
$arrS1R4[0,23] = 32887780.79329398
$arrS1R4[0,24] = 33019051.1094755
$arrS1R4[0,25] = 9079408.407006437
$arrS1R4[0,26] = 9247020.237221925
$arrS1R4[0,27] = 9446887.575191883
$arrS1R4[0,28] = 9637402.165375361
$arrS1R4[0,29] = 9882099.062900379
$arrS1R4[0,30] = 2634324.649063896
$arrS1R4[0,31] = 2885445.61736868
$arrS1R4[0,32] = 3304983.601109055
$arrS1R4[0,33] = 3911188.682747431
$arrS1R4[0,34] = 4681540.036266249
$arrS1R4[0,35] = 9105626.515739188
$arrS1R4[0,36] = 10368125.06706963
$arrS1R4[0,37] = 11839486.40330121
$arrS1R4[0,38] = 13507030.57311474
$arrS1R4[0,39] = 14784867.24470953
$arrS1R4[0,40] = 20653740.85302317
$arrS1R4[0,41] = 21342389.73276047
$arrS1R4[0,42] = 21830545.34489401
$arrS1R4[0,43] = 22092914.63860742
$arrS1R4[0,44] = 22170386.29645523
$arrS1R4[0,45] = 23341749.17863254
$arrS1R4[0,46] = 22768670.4645321
$arrS1R4[0,47] = 21644694.07133049
$arrS1R4[0,48] = 20143344.31843414
$arrS1R4[0,49] = 19026971.51524213
$arrS1R4[0,50] = 17369210.64525886
$arrS1R4[0,51] = 16660596.70824914
$arrS1R4[0,52] = 15975278.3139108
$arrS1R4[0,53] = 15431432.01748001
$arrS1R4[0,54] = 15107583.08799718
$arrS1R4[0,55] = 14516547.17774849
$arrS1R4[0,56] = 14718842.9073435
$arrS1R4[0,57] = 15316632.16414694
$arrS1R4[0,58] = 16244631.33828624
$arrS1R4[0,59] = 17353715.54407233
$arrS1R4[0,60] = 18452082.56303243
$arrS1R4[0,61] = 19734355.46247735
$arrS1R4[0,62] = 20938334.22034578
$arrS1R4[0,63] = 21984382.67961864
$arrS1R4[0,64] = 22776477.45966075
$arrS1R4[0,65] = 23334140.9353778
$arrS1R4[0,66] = 23781318.34468105
$arrS1R4[0,67] = 24081959.16740772
$arrS1R4[0,68] = 24159989.67191011
$arrS1R4[0,69] = 23911132.34547214
$arrS1R4[0,70] = 23215838.22245894
$arrS1R4[0,71] = 22556029.67007605
$arrS1R4[0,72] = 22075307.52441104
$arrS1R4[0,73] = 21789585.98889285
$arrS1R4[0,74] = 21489230.32532429
$arrS1R4[0,75] = 20962695.65789351
$arrS1R4[0,76] = 20666766.36915706
$arrS1R4[0,77] = 20452967.81678404
$arrS1R4[0,78] = 20359897.65675836
$arrS1R4[0,79] = 20404544.46762502
$arrS1R4[0,80] = 20415960.37196534
$arrS1R4[0,81] = 21035806.89924143
$arrS1R4[0,82] = 21779236.49194844
$arrS1R4[0,83] = 22585627.37432782
$arrS1R4[0,84] = 23405663.81587356
$arrS1R4[0,85] = 23600869.85625803
$arrS1R4[0,86] = 24241914.0630935
$arrS1R4[0,87] = 24739072.03732761
$arrS1R4[0,88] = 25111390.22202583
$arrS1R4[0,89] = 25386150.46503336
$ws1.Range("D4:CO4").Value = $arrS1R4

$arrS1R5 = New-Object 'object[,]' 1,90
$arrS1R5[0,0] = 2149741.752917891
$arrS1R5[0,1] = 1439679.414002341
$arrS1R5[0,2] = 1345818.204103041
$arrS1R5[0,3] = 1297234.071478467
$arrS1R5[0,4] = 1244834.330190728
$arrS1R5[0,5] = 1555389.286635582
$arrS1R5[0,6] = 1711377.353521905
$arrS1R5[0,7] = 1489426.94034161
$arrS1R5[0,8] = 1378420.280851263
$arrS1R5[0,9] = 1542529.424435778
$arrS1R5[0,10] = 1621887.367372605
$arrS1R5[0,11] = 2112184.618362186
$arrS1R5[0,12] = 16307854.80029605
$arrS1R5[0,13] = 16333811.00661036
$arrS1R5[0,14] = 16366767.76685695
$arrS1R5[0,15] = 54209877.12774869
$arrS1R5[0,16] = 54260358.68841778
$arrS1R5[0,17] = 54321571.15603209
$arrS1R5[0,18] = 54395315.13733557
$arrS1R5[0,19] = 54484202.46248683
$arrS1R5[0,20] = 57389938.31145861
$arrS1R5[0,21] = 57522856.98112737
$arrS1R5[0,22] = 57689253.9319164
$arrS1R5[0,23] = 57901201.50323297
$arrS1R5[0,24] = 58175259.8114974
$arrS1R5[0,25] = 16815757.13610247
$arrS1R5[0,26] = 17286365.34550833
$arrS1R5[0,27] = 17904375.76482018
$arrS1R5[0,28] = 18711087.27449148
$arrS1R5[0,29] = 19752638.32282598
$arrS1R5[0,30] = 7100985.338434661
$arrS1R5[0,31] = 8558628.694068663
$arrS1R5[0,32] = 10343223.17292801
$arrS1R5[0,33] = 12736317.07990005
$arrS1R5[0,34] = 15543797.39676337
$arrS1R5[0,35] = 24623415.36617304
$arrS1R5[0,36] = 28103846.88238889
$arrS1R5[0,37] = 31788746.97210365
$arrS1R5[0,38] = 35538735.84028139
$arrS1R5[0,39] = 39182083.57401909
$arrS1R5[0,40] = 50730210.09694149
$arrS1R5[0,41] = 53582963.06524845
$arrS1R5[0,42] = 55768524.12185865
$arrS1R5[0,43] = 57152310.57107791
$arrS1R5[0,44] = 57657577.4055163
$arrS1R5[0,45] = 59674793.83496889
$arrS1R5[0,46] = 58466664.70073942
$arrS1R5[0,47] = 56555590.35647823
$arrS1R5[0,48] = 54114422.25143787
$arrS1R5[0,49] = 51349754.39270806
$arrS1R5[0,50] = 47088943.78502554
$arrS1R5[0,51] = 44351865.63733242
$arrS1R5[0,52] = 41954718.08378091
$arrS1R5[0,53] = 40081560.01682471
$arrS1R5[0,54] = 38872385.50926251
$arrS1R5[0,55] = 37640546.0234106
$arrS1R5[0,56] = 37952029.72608801
$arrS1R5[0,57] = 39002643.71591939
$arrS1R5[0,58] = 40712331.19998205
$arrS1R5[0,59] = 42967377.25583482
$arrS1R5[0,60] = 45527530.78721425
$arrS1R5[0,61] = 48457143.71103125
$arrS1R5[0,62] = 51500195.74780367
$arrS1R5[0,63] = 54508456.68066286
$arrS1R5[0,64] = 57343204.35894439
$arrS1R5[0,65] = 59906786.6344051
$arrS1R5[0,66] = 62052452.1056831
$arrS1R5[0,67] = 63734008.00617307
$arrS1R5[0,68] = 64911880.09468625
$arrS1R5[0,69] = 65578272.31615769
$arrS1R5[0,70] = 65537198.36493602
$arrS1R5[0,71] = 65280738.10258959
$arrS1R5[0,72] = 64667631.89649016
$arrS1R5[0,73] = 63795663.61114269
$arrS1R5[0,74] = 62774651.21854672
$arrS1R5[0,75] = 61327560.50123087
$arrS1R5[0,76] = 60346309.92846997
$arrS1R5[0,77] = 59539220.10646011
$arrS1R5[0,78] = 58988733.0391894
$arrS1R5[0,79] = 58755790.58340008
$arrS1R5[0,80] = 58206316.45132345
$arrS1R5[0,81] = 58691978.72576661
$arrS1R5[0,82] = 59528162.71852187
$arrS1R5[0,83] = 60678712.23020957
$arrS1R5[0,84] = 62089476.97916714
$arrS1R5[0,85] = 62640841.93142905
$arrS1R5[0,86] = 64362397.59358234
$arrS1R5[0,87] = 66123977.81259896
$arrS1R5[0,88] = 67850172.64941625
$arrS1R5[0,89] = 69472419.29161821
$ws1.Range("D5:CO5").Value = $arrS1R5

$arrS2R2 = New-Object 'object[,]' 1,90
$arrS2R2[0,0] = 1072.353295692879
$arrS2R2[0,1] = 2469.945846677267
$arrS2R2[0,2] = 3866.852345457151
$arrS2R2[0,3] = 5194.288700394052
$arrS2R2[0,4] = 6974.68540242248
$arrS2R2[0,5] = 8259.113854798166
$arrS2R2[0,6] = 10583.19521345397
$arrS2R2[0,7] = 12893.65018086489
$arrS2R2[0,8] = 14949.47695446425
$arrS2R2[0,9] = 16844.38229612978
$arrS2R2[0,10] = 20897.987938082
$arrS2R2[0,11] = 23358.57051525402
$arrS2R2[0,12] = 28290.22427334399
$arrS2R2[0,13] = 33612.63818200492
$arrS2R2[0,14] = 39331.25078999881
$arrS2R2[0,15] = 51834.42073356161
$arrS2R2[0,16] = 65221.81400187083
$arrS2R2[0,17] = 79503.03179033326
$arrS2R2[0,18] = 94694.920168146
$arrS2R2[0,19] = 110810.8861707333
$arrS2R2[0,20] = 127005.5210793848
$arrS2R2[0,21] = 144089.1118818538
$arrS2R2[0,22] = 162069.0130647209
$arrS2R2[0,23] = 180963.993300658
$arrS2R2[0,24] = 200799.9644805966
$arrS2R2[0,25] = 210074.8410006412
$arrS2R2[0,26] = 219607.7458095777
$arrS2R2[0,27] = 229471.7635563455
$arrS2R2[0,28] = 239656.3183578329
$arrS2R2[0,29] = 250215.7436503791
$arrS2R2[0,30] = 256987.5923063907
$arrS2R2[0,31] = 263822.9133434905
$arrS2R2[0,32] = 270989.5448516527
$arrS2R2[0,33] = 278748.8710352551
$arrS2R2[0,34] = 287085.3328535687
$arrS2R2[0,35] = 298129.8913284243
$arrS2R2[0,36] = 310571.6654617666
$arrS2R2[0,37] = 324950.3441595362
$arrS2R2[0,38] = 341828.5144372305
$arrS2R2[0,39] = 360497.9725315598
$arrS2R2[0,40] = 383479.0492214851
$arrS2R2[0,41] = 407450.4194562975
$arrS2R2[0,42] = 432253.4167787464
$arrS2R2[0,43] = 457736.2371679994
$arrS2R2[0,44] = 483909.0763073803
$arrS2R2[0,45] = 511323.4017162066
$arrS2R2[0,46] = 538462.3381156095
$arrS2R2[0,47] = 563840.2019318022
$arrS2R2[0,48] = 586357.9286650572
$arrS2R2[0,49] = 607456.567281703
$arrS2R2[0,50] = 627673.2433159193
$arrS2R2[0,51] = 647754.2513769221
$arrS2R2[0,52] = 667194.7128140718
$arrS2R2[0,53] = 685655.6771343343
$arrS2R2[0,54] = 703795.8901769413
$arrS2R2[0,55] = 721759.7266527959
$arrS2R2[0,56] = 740575.3664785388
$arrS2R2[0,57] = 760699.2062728935
$arrS2R2[0,58] = 782618.9388414292
$arrS2R2[0,59] = 806501.9613528876
$arrS2R2[0,60] = 832125.9021683608
$arrS2R2[0,61] = 860080.4699557148
$arrS2R2[0,62] = 890435.8725754713
$arrS2R2[0,63] = 923047.4570258493
$arrS2R2[0,64] = 957194.3439210402
$arrS2R2[0,65] = 992085.9166887921
$arrS2R2[0,66] = 1027761.080027909
$arrS2R2[0,67] = 1064073.29078854
$arrS2R2[0,68] = 1100613.668803032
$arrS2R2[0,69] = 1136525.238920912
$arrS2R2[0,70] = 1170735.496794471
$arrS2R2[0,71] = 1203264.341029952
$arrS2R2[0,72] = 1234608.112449533
$arrS2R2[0,73] = 1265325.111712801
$arrS2R2[0,74] = 1295638.19311247
$arrS2R2[0,75] = 1325300.064561385
$arrS2R2[0,76] = 1354547.031753223
$arrS2R2[0,77] = 1383723.269230985
$arrS2R2[0,78] = 1413330.111766741
$arrS2R2[0,79] = 1443931.641588225
$arrS2R2[0,80] = 1475747.30180299
$arrS2R2[0,81] = 1509337.226238863
$arrS2R2[0,82] = 1544862.434385323
$arrS2R2[0,83] = 1582422.775433045
$arrS2R2[0,84] = 1622143.564560264
$arrS2R2[0,85] = 1663590.387295622
$arrS2R2[0,86] = 1706846.550163591
$arrS2R2[0,87] = 1751421.460785394
$arrS2R2[0,88] = 1796869.505557852
$arrS2R2[0,89] = 1842917.295855446
$ws2.Range("D2:CO2").Value = $arrS2R2

$arrS2R3 = New-Object 'object[,]' 1,90
$arrS2R3[0,0] = 70.9098667458114
$arrS2R3[0,1] = 163.6056466423233
$arrS2R3[0,2] = 256.13708219816
$arrS2R3[0,3] = 344.001556144818
$arrS2R3[0,4] = 462.3864842671475
$arrS2R3[0,5] = 547.4485610075349
$arrS2R3[0,6] = 702.2250829589436
$arrS2R3[0,7] = 855.9341932330033
$arrS2R3[0,8] = 992.5869185861493
$arrS2R3[0,9] = 1118.354362807883
$arrS2R3[0,10] = 1389.019764050263
$arrS2R3[0,11] = 1552.830135403919
$arrS2R3[0,12] = 1882.485974792283
$arrS2R3[0,13] = 2238.359400880697
$arrS2R3[0,14] = 2620.818856303182
$arrS2R3[0,15] = 3458.775713157263
$arrS2R3[0,16] = 4356.073065713931
$arrS2R3[0,17] = 5313.360751274489
$arrS2R3[0,18] = 6331.779280754066
$arrS2R3[0,19] = 7412.236246926958
$arrS2R3[0,20] = 8497.962603417407
$arrS2R3[0,21] = 9643.396450717664
$arrS2R3[0,22] = 10849.03534844354
$arrS2R3[0,23] = 12116.14936929955
$arrS2R3[0,24] = 13446.49203721808
$arrS2R3[0,25] = 14067.67651094418
$arrS2R3[0,26] = 14706.21099335574
$arrS2R3[0,27] = 15367.04188579259
$arrS2R3[0,28] = 16049.44956522021
$arrS2R3[0,29] = 16757.10852678629
$arrS2R3[0,30] = 17210.44845794315
$arrS2R3[0,31] = 17668.05607687437
$arrS2R3[0,32] = 18148.06275365501
$arrS2R3[0,33] = 18668.15752904375
$arrS2R3[0,34] = 19227.27634330161
$arrS2R3[0,35] = 19968.70329232487
$arrS2R3[0,36] = 20804.63184844841
$arrS2R3[0,37] = 21771.60503917957
$arrS2R3[0,38] = 22907.71977261407
$arrS2R3[0,39] = 24165.02139973271
$arrS2R3[0,40] = 25712.3237947604
$arrS2R3[0,41] = 27326.51389873635
$arrS2R3[0,42] = 28996.86140552456
$arrS2R3[0,43] = 30713.10363174919
$arrS2R3[0,44] = 32475.95007971457
$arrS2R3[0,45] = 34322.33910369307
$arrS2R3[0,46] = 36149.97189693576
$arrS2R3[0,47] = 37858.25160887317
$arrS2R3[0,48] = 39372.75936262778
$arrS2R3[0,49] = 40791.09041024763
$arrS2R3[0,50] = 42150.12796895189
$arrS2R3[0,51] = 43499.91259269226
$arrS2R3[0,52] = 44806.24474399709
$arrS2R3[0,53] = 46046.15734831411
$arrS2R3[0,54] = 47264.24826003316
$arrS2R3[0,55] = 48470.64504101018
$arrS2R3[0,56] = 49734.60209106243
$arrS2R3[0,57] = 51087.00957413706
$arrS2R3[0,58] = 52560.87918651651
$arrS2R3[0,59] = 54167.53490338743
$arrS2R3[0,60] = 55892.11042042025
$arrS2R3[0,61] = 57774.33756610031
$arrS2R3[0,62] = 59818.97077854156
$arrS2R3[0,63] = 62016.21691049764
$arrS2R3[0,64] = 64317.26136841045
$arrS2R3[0,65] = 66668.66882269818
$arrS2R3[0,66] = 69072.98017688619
$arrS2R3[0,67] = 71520.28155557564
$arrS2R3[0,68] = 73982.89413233618
$arrS2R3[0,69] = 76402.79346311325
$arrS2R3[0,70] = 78707.48172116927
$arrS2R3[0,71] = 80898.19032726111
$arrS2R3[0,72] = 83008.53817687643
$arrS2R3[0,73] = 85076.33815950758
$arrS2R3[0,74] = 87116.68564716175
$arrS2R3[0,75] = 89112.9463286676
$arrS2R3[0,76] = 91081.0068765412
$arrS2R3[0,77] = 93044.17203162142
$arrS2R3[0,78] = 95036.38804759647
$arrS2R3[0,79] = 97095.84912113316
$arrS2R3[0,80] = 99237.5489305724
$arrS2R3[0,81] = 101499.2665227509
$arrS2R3[0,82] = 103891.8985007469
$arrS2R3[0,83] = 106422.2004417968
$arrS2R3[0,84] = 109098.6545669464
$arrS2R3[0,85] = 111892.0618217478
$arrS2R3[0,86] = 114807.838790714
$arrS2R3[0,87] = 117812.7615434047
$arrS2R3[0,88] = 120876.6533492622
$arrS2R3[0,89] = 123981.0023022256
$ws2.Range("D3:CO3").Value = $arrS2R3

$arrS2R4 = New-Object 'object[,]' 1,90
$arrS2R4[0,0] = 1624090.537329962
$arrS2R4[0,1] = 3082913.836292786
$arrS2R4[0,2] = 4356450.785156812
$arrS2R4[0,3] = 5833709.157689445
$arrS2R4[0,4] = 7457067.915503096
$arrS2R4[0,5] = 8999940.812731894
$arrS2R4[0,6] = 10677359.86967593
$arrS2R4[0,7] = 12215980.65673972
$arrS2R4[0,8] = 13744287.09156763
$arrS2R4[0,9] = 15782708.16079887
$arrS2R4[0,10] = 17717116.39726498
$arrS2R4[0,11] = 19633269.52412627
$arrS2R4[0,12] = 29564847.36441307
$arrS2R4[0,13] = 39547234.64284137
$arrS2R4[0,14] = 49569665.69633654
$arrS2R4[0,15] = 80849446.82505143
$arrS2R4[0,16] = 112193824.1416991
$arrS2R4[0,17] = 143613811.0749336
$arrS2R4[0,18] = 175119339.1718884
$arrS2R4[0,19] = 206718902.4534993
$arrS2R4[0,20] = 239268863.8547618
$arrS2R4[0,21] = 271929775.4641445
$arrS2R4[0,22] = 304699086.8125048
$arrS2R4[0,23] = 337586867.6057988
$arrS2R4[0,24] = 370605918.7152743
$arrS2R4[0,25] = 379685327.1222807
$arrS2R4[0,26] = 388932347.3595026
$arrS2R4[0,27] = 398379234.9346945
$arrS2R4[0,28] = 408016637.1000698
$arrS2R4[0,29] = 417898736.1629702
$arrS2R4[0,30] = 420533060.8120341
$arrS2R4[0,31] = 423418506.4294028
$arrS2R4[0,32] = 426723490.0305118
$arrS2R4[0,33] = 430634678.7132592
$arrS2R4[0,34] = 435316218.7495255
$arrS2R4[0,35] = 444421845.2652647
$arrS2R4[0,36] = 454789970.3323343
$arrS2R4[0,37] = 466629456.7356355
$arrS2R4[0,38] = 480136487.3087503
$arrS2R4[0,39] = 494921354.5534598
$arrS2R4[0,40] = 515575095.406483
$arrS2R4[0,41] = 536917485.1392435
$arrS2R4[0,42] = 558748030.4841375
$arrS2R4[0,43] = 580840945.1227449
$arrS2R4[0,44] = 603011331.4192002
$arrS2R4[0,45] = 626353080.5978327
$arrS2R4[0,46] = 649121751.0623648
$arrS2R4[0,47] = 670766445.1336954
$arrS2R4[0,48] = 690909789.4521295
$arrS2R4[0,49] = 709936760.9673716
$arrS2R4[0,50] = 727305971.6126305
$arrS2R4[0,51] = 743966568.3208796
$arrS2R4[0,52] = 759941846.6347904
$arrS2R4[0,53] = 775373278.6522704
$arrS2R4[0,54] = 790480861.7402676
$arrS2R4[0,55] = 804997408.9180161
$arrS2R4[0,56] = 819716251.8253596
$arrS2R4[0,57] = 835032883.9895065
$arrS2R4[0,58] = 851277515.3277928
$arrS2R4[0,59] = 868631230.871865
$arrS2R4[0,60] = 887083313.4348974
$arrS2R4[0,61] = 906817668.8973747
$arrS2R4[0,62] = 927756003.1177205
$arrS2R4[0,63] = 949740385.7973391
$arrS2R4[0,64] = 972516863.2569999
$arrS2R4[0,65] = 995851004.1923777
$arrS2R4[0,66] = 1019632322.537059
$arrS2R4[0,67] = 1043714281.704466
$arrS2R4[0,68] = 1067874271.376377
$arrS2R4[0,69] = 1091785403.721849
$arrS2R4[0,70] = 1115001241.944308
$arrS2R4[0,71] = 1137557271.614384
$arrS2R4[0,72] = 1159632579.138795
$arrS2R4[0,73] = 1181422165.127687
$arrS2R4[0,74] = 1202911395.453012
$arrS2R4[0,75] = 1223874091.110905
$arrS2R4[0,76] = 1244540857.480062
$arrS2R4[0,77] = 1264993825.296846
$arrS2R4[0,78] = 1285353722.953605
$arrS2R4[0,79] = 1305758267.42123
$arrS2R4[0,80] = 1326174227.793195
$arrS2R4[0,81] = 1347210034.692437
$arrS2R4[0,82] = 1368989271.184385
$arrS2R4[0,83] = 1391574898.558713
$arrS2R4[0,84] = 1414980562.374587
$arrS2R4[0,85] = 1438581432.230844
$arrS2R4[0,86] = 1462823346.293938
$arrS2R4[0,87] = 1487562418.331265
$arrS2R4[0,88] = 1512673808.553291
$arrS2R4[0,89] = 1538059959.018325
$ws2.Range("D4:CO4").Value = $arrS2R4

$arrS2R5 = New-Object 'object[,]' 1,90
$arrS2R5[0,0] = 2149741.752917891
$arrS2R5[0,1] = 3589421.166920232
$arrS2R5[0,2] = 4935239.371023273
$arrS2R5[0,3] = 6232473.44250174
$arrS2R5[0,4] = 7477307.772692468
$arrS2R5[0,5] = 9032697.059328049
$arrS2R5[0,6] = 10744074.41284996
$arrS2R5[0,7] = 12233501.35319157
$arrS2R5[0,8] = 13611921.63404283
$arrS2R5[0,9] = 15154451.05847861
$arrS2R5[0,10] = 16776338.42585121
$arrS2R5[0,11] = 18888523.0442134
$arrS2R5[0,12] = 35196377.84450945
$arrS2R5[0,13] = 51530188.85111981
$arrS2R5[0,14] = 67896956.61797675
$arrS2R5[0,15] = 122106833.7457255
$arrS2R5[0,16] = 176367192.4341432
$arrS2R5[0,17] = 230688763.5901753
$arrS2R5[0,18] = 285084078.7275109
$arrS2R5[0,19] = 339568281.1899977
$arrS2R5[0,20] = 396958219.5014563
$arrS2R5[0,21] = 454481076.4825837
$arrS2R5[0,22] = 512170330.4145001
$arrS2R5[0,23] = 570071531.9177331
$arrS2R5[0,24] = 628246791.7292305
$arrS2R5[0,25] = 645062548.865333
$arrS2R5[0,26] = 662348914.2108413
$arrS2R5[0,27] = 680253289.9756615
$arrS2R5[0,28] = 698964377.2501529
$arrS2R5[0,29] = 718717015.572979
$arrS2R5[0,30] = 725818000.9114137
$arrS2R5[0,31] = 734376629.6054823
$arrS2R5[0,32] = 744719852.7784103
$arrS2R5[0,33] = 757456169.8583103
$arrS2R5[0,34] = 772999967.2550737
$arrS2R5[0,35] = 797623382.6212467
$arrS2R5[0,36] = 825727229.5036355
$arrS2R5[0,37] = 857515976.4757391
$arrS2R5[0,38] = 893054712.3160205
$arrS2R5[0,39] = 932236795.8900396
$arrS2R5[0,40] = 982967005.986981
$arrS2R5[0,41] = 1036549969.05223
$arrS2R5[0,42] = 1092318493.174088
$arrS2R5[0,43] = 1149470803.745166
$arrS2R5[0,44] = 1207128381.150682
$arrS2R5[0,45] = 1266803174.985651
$arrS2R5[0,46] = 1325269839.686391
$arrS2R5[0,47] = 1381825430.042869
$arrS2R5[0,48] = 1435939852.294307
$arrS2R5[0,49] = 1487289606.687015
$arrS2R5[0,50] = 1534378550.47204
$arrS2R5[0,51] = 1578730416.109373
$arrS2R5[0,52] = 1620685134.193154
$arrS2R5[0,53] = 1660766694.209979
$arrS2R5[0,54] = 1699639079.719241
$arrS2R5[0,55] = 1737279625.742652
$arrS2R5[0,56] = 1775231655.46874
$arrS2R5[0,57] = 1814234299.184659
$arrS2R5[0,58] = 1854946630.384641
$arrS2R5[0,59] = 1897914007.640476
$arrS2R5[0,60] = 1943441538.427691
$arrS2R5[0,61] = 1991898682.138722
$arrS2R5[0,62] = 2043398877.886525
$arrS2R5[0,63] = 2097907334.567188
$arrS2R5[0,64] = 2155250538.926133
$arrS2R5[0,65] = 2215157325.560538
$arrS2R5[0,66] = 2277209777.666221
$arrS2R5[0,67] = 2340943785.672394
$arrS2R5[0,68] = 2405855665.76708
$arrS2R5[0,69] = 2471433938.083238
$arrS2R5[0,70] = 2536971136.448174
$arrS2R5[0,71] = 2602251874.550764
$arrS2R5[0,72] = 2666919506.447254
$arrS2R5[0,73] = 2730715170.058396
$arrS2R5[0,74] = 2793489821.276943
$arrS2R5[0,75] = 2854817381.778174
$arrS2R5[0,76] = 2915163691.706644
$arrS2R5[0,77] = 2974702911.813104
$arrS2R5[0,78] = 3033691644.852293
$arrS2R5[0,79] = 3092447435.435694
$arrS2R5[0,80] = 3150653751.887017
$arrS2R5[0,81] = 3209345730.612784
$arrS2R5[0,82] = 3268873893.331306
$arrS2R5[0,83] = 3329552605.561515
$arrS2R5[0,84] = 3391642082.540682
$arrS2R5[0,85] = 3454282924.472111
$arrS2R5[0,86] = 3518645322.065693
$arrS2R5[0,87] = 3584769299.878293
$arrS2R5[0,88] = 3652619472.527709
$arrS2R5[0,89] = 3722091891.819327
$ws2.Range("D5:CO5").Value = $arrS2R5
